# Generate Report for Handoff
#
# The localization status report is regenerated for handoff: the overall
# status moves from "Handed back: in sync with en-US" to "Ready for
# handoff", the associated timestamps advance, and the now-shorter status
# text lets the status/date columns shrink on the Overview, zh-cn and
# de-de sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
$newColWidth = 17.2159881591797

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-10-21 04:32:59"

$overview.Columns.Item(5).ColumnWidth = $newColWidth
$overview.Columns.Item(6).ColumnWidth = $newColWidth

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-10-21 04:32:48"
$zhcn.Columns.Item(3).ColumnWidth = $newColWidth

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-10-21 04:32:59"
$dede.Columns.Item(3).ColumnWidth = $newColWidth
